# Update "Countries & provincias Spain" worksheet:
#  - refresh the "last updated" timestamp
#  - refresh case counts for several provinces
#  - two provinces (Gran Canaria/Huesca and Arroyo de la Luz/La Gomera)
#    swap position because their totals crossed over during the refresh

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 10 de Abril de 2020 a las 20:22"

# Row 32: Tenerife
$ws.Cells.Item(32, 1).Value = "Tenerife"
$ws.Cells.Item(32, 2).Value = 1238
$ws.Cells.Item(32, 3).Value = 281
$ws.Cells.Item(32, 4).Value = 894
$ws.Cells.Item(32, 5).Value = 63

# Row 47: was Huesca, now Gran Canaria (overtook Huesca in total cases)
$ws.Cells.Item(47, 1).Value = "Gran Canaria"
$ws.Cells.Item(47, 2).Value = 480
$ws.Cells.Item(47, 3).Value = 97
$ws.Cells.Item(47, 4).Value = 356
$ws.Cells.Item(47, 5).Value = 27

# Row 48: was Gran Canaria, now Huesca
$ws.Cells.Item(48, 1).Value = "Huesca"
$ws.Cells.Item(48, 2).Value = 478
$ws.Cells.Item(48, 3).Value = 83
$ws.Cells.Item(48, 4).Value = 331
$ws.Cells.Item(48, 5).Value = 64

# Row 56: La Palma
$ws.Cells.Item(56, 1).Value = "La Palma"
$ws.Cells.Item(56, 2).Value = 71
$ws.Cells.Item(56, 3).Value = 9
$ws.Cells.Item(56, 4).Value = 59
$ws.Cells.Item(56, 5).Value = 3

# Row 57: Lanzarote
$ws.Cells.Item(57, 1).Value = "Lanzarote"
$ws.Cells.Item(57, 2).Value = 67
$ws.Cells.Item(57, 3).Value = 12
$ws.Cells.Item(57, 4).Value = 53
$ws.Cells.Item(57, 5).Value = 2

# Row 59: Fuerteventura
$ws.Cells.Item(59, 1).Value = "Fuerteventura"
$ws.Cells.Item(59, 2).Value = 24
$ws.Cells.Item(59, 3).Value = 7
$ws.Cells.Item(59, 4).Value = 17
$ws.Cells.Item(59, 5).Value = 0

# Row 62: was La Gomera, now Arroyo de la Luz
$ws.Cells.Item(62, 1).Value = "Arroyo de la Luz"
$ws.Cells.Item(62, 2).Value = 7
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(62, 4).Value = 7
$ws.Cells.Item(62, 5).Value = 0

# Row 63: was Arroyo de la Luz, now La Gomera
$ws.Cells.Item(63, 1).Value = "La Gomera"
$ws.Cells.Item(63, 2).Value = 7
$ws.Cells.Item(63, 3).Value = 2
$ws.Cells.Item(63, 4).Value = 5
$ws.Cells.Item(63, 5).Value = 0

# Row 64: El Hierro
$ws.Cells.Item(64, 1).Value = "El Hierro"
$ws.Cells.Item(64, 2).Value = 1
$ws.Cells.Item(64, 3).Value = 1
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
